# Auto-generated Excel COM-interop script
# Updates countries & provincias Spain dataset to reflect the 10 Apr 2020 06:52 snapshot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 10 de Abril de 2020 a las 06:52
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 10 de Abril de 2020 a las 06:52"

# Row 36: Pakistan
$ws.Cells.Item(36, 1).Value2 = "Pakistan"
$ws.Cells.Item(36, 2).Value2 = 4601
$ws.Cells.Item(36, 3).Value2 = 112
$ws.Cells.Item(36, 4).Value2 = 572
$ws.Cells.Item(36, 5).Value2 = 3963
$ws.Cells.Item(36, 6).Value2 = 31
$ws.Cells.Item(36, 7).Value2 = 1
$ws.Cells.Item(36, 8).Value2 = 66

# Row 53: Singapur
$ws.Cells.Item(53, 1).Value2 = "Singapur"
$ws.Cells.Item(53, 2).Value2 = 1910
$ws.Cells.Item(53, 3).Value2 = 0
$ws.Cells.Item(53, 4).Value2 = 460
$ws.Cells.Item(53, 5).Value2 = 1444
$ws.Cells.Item(53, 6).Value2 = 29
$ws.Cells.Item(53, 7).Value2 = 0
$ws.Cells.Item(53, 8).Value2 = 6

# Row 67: Lituania
$ws.Cells.Item(67, 1).Value2 = "Lituania"
$ws.Cells.Item(67, 2).Value2 = 999
$ws.Cells.Item(67, 3).Value2 = 44
$ws.Cells.Item(67, 4).Value2 = 8
$ws.Cells.Item(67, 5).Value2 = 975
$ws.Cells.Item(67, 6).Value2 = 21
$ws.Cells.Item(67, 7).Value2 = 0
$ws.Cells.Item(67, 8).Value2 = 16

# Row 68: Hungria
$ws.Cells.Item(68, 1).Value2 = "Hungria"
$ws.Cells.Item(68, 2).Value2 = 980
$ws.Cells.Item(68, 3).Value2 = 0
$ws.Cells.Item(68, 4).Value2 = 96
$ws.Cells.Item(68, 5).Value2 = 818
$ws.Cells.Item(68, 6).Value2 = 17
$ws.Cells.Item(68, 7).Value2 = 0
$ws.Cells.Item(68, 8).Value2 = 66

# Row 69: Hong Kong
$ws.Cells.Item(69, 1).Value2 = "Hong Kong"
$ws.Cells.Item(69, 2).Value2 = 974
$ws.Cells.Item(69, 3).Value2 = 0
$ws.Cells.Item(69, 4).Value2 = 293
$ws.Cells.Item(69, 5).Value2 = 677
$ws.Cells.Item(69, 6).Value2 = 14
$ws.Cells.Item(69, 7).Value2 = 0
$ws.Cells.Item(69, 8).Value2 = 4

# Row 96: Honduras
$ws.Cells.Item(96, 1).Value2 = "Honduras"
$ws.Cells.Item(96, 2).Value2 = 382
$ws.Cells.Item(96, 3).Value2 = 39
$ws.Cells.Item(96, 4).Value2 = 7
$ws.Cells.Item(96, 5).Value2 = 352
$ws.Cells.Item(96, 6).Value2 = 10
$ws.Cells.Item(96, 7).Value2 = 0
$ws.Cells.Item(96, 8).Value2 = 23

# Row 128: Camboya
$ws.Cells.Item(128, 1).Value2 = "Camboya"
$ws.Cells.Item(128, 2).Value2 = 119
$ws.Cells.Item(128, 3).Value2 = 0
$ws.Cells.Item(128, 4).Value2 = 72
$ws.Cells.Item(128, 5).Value2 = 47
$ws.Cells.Item(128, 6).Value2 = 1
$ws.Cells.Item(128, 7).Value2 = 0
$ws.Cells.Item(128, 8).Value2 = 0

# Row 145: San Martin (Parte Holandesa)
$ws.Cells.Item(145, 1).Value2 = "San Martin (Parte Holandesa)"
$ws.Cells.Item(145, 2).Value2 = 50
$ws.Cells.Item(145, 3).Value2 = 7
$ws.Cells.Item(145, 4).Value2 = 3
$ws.Cells.Item(145, 5).Value2 = 39
$ws.Cells.Item(145, 6).Value2 = 2
$ws.Cells.Item(145, 7).Value2 = 2
$ws.Cells.Item(145, 8).Value2 = 8

# Row 146: Bermudas
$ws.Cells.Item(146, 1).Value2 = "Bermudas"
$ws.Cells.Item(146, 2).Value2 = 48
$ws.Cells.Item(146, 3).Value2 = 0
$ws.Cells.Item(146, 4).Value2 = 25
$ws.Cells.Item(146, 5).Value2 = 19
$ws.Cells.Item(146, 6).Value2 = 2
$ws.Cells.Item(146, 7).Value2 = 0
$ws.Cells.Item(146, 8).Value2 = 4

# Row 147: Islas Caimanes
$ws.Cells.Item(147, 1).Value2 = "Islas Caimanes"
$ws.Cells.Item(147, 2).Value2 = 45
$ws.Cells.Item(147, 3).Value2 = 0
$ws.Cells.Item(147, 4).Value2 = 6
$ws.Cells.Item(147, 5).Value2 = 38
$ws.Cells.Item(147, 6).Value2 = 0
$ws.Cells.Item(147, 7).Value2 = 0
$ws.Cells.Item(147, 8).Value2 = 1

# Row 148: Macao
$ws.Cells.Item(148, 1).Value2 = "Macao"
$ws.Cells.Item(148, 2).Value2 = 45
$ws.Cells.Item(148, 3).Value2 = 0
$ws.Cells.Item(148, 4).Value2 = 10
$ws.Cells.Item(148, 5).Value2 = 35
$ws.Cells.Item(148, 6).Value2 = 1
$ws.Cells.Item(148, 7).Value2 = 0
$ws.Cells.Item(148, 8).Value2 = 0

# Row 149: Gabon
$ws.Cells.Item(149, 1).Value2 = "Gabon"
$ws.Cells.Item(149, 2).Value2 = 44
$ws.Cells.Item(149, 3).Value2 = 0
$ws.Cells.Item(149, 4).Value2 = 1
$ws.Cells.Item(149, 5).Value2 = 42
$ws.Cells.Item(149, 6).Value2 = 0
$ws.Cells.Item(149, 7).Value2 = 0
$ws.Cells.Item(149, 8).Value2 = 1

# Row 172: Laos
$ws.Cells.Item(172, 1).Value2 = "Laos"
$ws.Cells.Item(172, 2).Value2 = 16
$ws.Cells.Item(172, 3).Value2 = 0
$ws.Cells.Item(172, 4).Value2 = 0
$ws.Cells.Item(172, 5).Value2 = 16
$ws.Cells.Item(172, 6).Value2 = 0
$ws.Cells.Item(172, 7).Value2 = 0
$ws.Cells.Item(172, 8).Value2 = 0

# Row 173: Fiyi
$ws.Cells.Item(173, 1).Value2 = "Fiyi"
$ws.Cells.Item(173, 2).Value2 = 16
$ws.Cells.Item(173, 3).Value2 = 1
$ws.Cells.Item(173, 4).Value2 = 0
$ws.Cells.Item(173, 5).Value2 = 16
$ws.Cells.Item(173, 6).Value2 = 0
$ws.Cells.Item(173, 7).Value2 = 0
$ws.Cells.Item(173, 8).Value2 = 0

# Row 209: Burundi
$ws.Cells.Item(209, 1).Value2 = "Burundi"
$ws.Cells.Item(209, 2).Value2 = 3
$ws.Cells.Item(209, 3).Value2 = 0
$ws.Cells.Item(209, 4).Value2 = 0
$ws.Cells.Item(209, 5).Value2 = 3
$ws.Cells.Item(209, 6).Value2 = 0
$ws.Cells.Item(209, 7).Value2 = 0
$ws.Cells.Item(209, 8).Value2 = 0

# Row 210: Islas Virgenes Britanicas
$ws.Cells.Item(210, 1).Value2 = "Islas Virgenes Britanicas"
$ws.Cells.Item(210, 2).Value2 = 3
$ws.Cells.Item(210, 3).Value2 = 0
$ws.Cells.Item(210, 4).Value2 = 0
$ws.Cells.Item(210, 5).Value2 = 3
$ws.Cells.Item(210, 6).Value2 = 0
$ws.Cells.Item(210, 7).Value2 = 0
$ws.Cells.Item(210, 8).Value2 = 0

# Row 211: Sudan del Sur
$ws.Cells.Item(211, 1).Value2 = "Sudan del Sur"
$ws.Cells.Item(211, 2).Value2 = 3
$ws.Cells.Item(211, 3).Value2 = 0
$ws.Cells.Item(211, 4).Value2 = 0
$ws.Cells.Item(211, 5).Value2 = 3
$ws.Cells.Item(211, 6).Value2 = 0
$ws.Cells.Item(211, 7).Value2 = 0
$ws.Cells.Item(211, 8).Value2 = 0

# Row 214: Yemen
$ws.Cells.Item(214, 1).Value2 = "Yemen"
$ws.Cells.Item(214, 2).Value2 = 1
$ws.Cells.Item(214, 3).Value2 = 1
$ws.Cells.Item(214, 4).Value2 = 0
$ws.Cells.Item(214, 5).Value2 = 1
$ws.Cells.Item(214, 6).Value2 = 0
$ws.Cells.Item(214, 7).Value2 = 0
$ws.Cells.Item(214, 8).Value2 = 0

# Row 215: Timor Oriental
$ws.Cells.Item(215, 1).Value2 = "Timor Oriental"
$ws.Cells.Item(215, 2).Value2 = 1
$ws.Cells.Item(215, 3).Value2 = 0
$ws.Cells.Item(215, 4).Value2 = 0
$ws.Cells.Item(215, 5).Value2 = 1
$ws.Cells.Item(215, 6).Value2 = 0
$ws.Cells.Item(215, 7).Value2 = 0
$ws.Cells.Item(215, 8).Value2 = 0

# Row 216: San Pedro y Miquelon
$ws.Cells.Item(216, 1).Value2 = "San Pedro y Miquelon"
$ws.Cells.Item(216, 2).Value2 = 1
$ws.Cells.Item(216, 3).Value2 = 0
$ws.Cells.Item(216, 4).Value2 = 0
$ws.Cells.Item(216, 5).Value2 = 1
$ws.Cells.Item(216, 6).Value2 = 0
$ws.Cells.Item(216, 7).Value2 = 0
$ws.Cells.Item(216, 8).Value2 = 0
